$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from E1 to F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timeValues = @(
    "2021-10-05 10:51:07.856377",
    "2021-10-05 10:51:07.856389",
    "2021-10-05 10:51:07.856392",
    "2021-10-05 10:51:07.856395",
    "2021-10-05 10:51:07.856398",
    "2021-10-05 10:51:07.856400",
    "2021-10-05 10:51:07.856403",
    "2021-10-05 10:51:07.856406",
    "2021-10-05 10:51:07.856408",
    "2021-10-05 10:51:07.856411",
    "2021-10-05 10:51:07.856414",
    "2021-10-05 10:51:07.856416",
    "2021-10-05 10:51:07.856419",
    "2021-10-05 10:51:07.856421",
    "2021-10-05 10:51:07.856424",
    "2021-10-05 10:51:07.856426",
    "2021-10-05 10:51:07.856429",
    "2021-10-05 10:51:07.856432",
    "2021-10-05 10:51:07.856434",
    "2021-10-05 10:51:07.856437",
    "2021-10-05 10:51:07.856439",
    "2021-10-05 10:51:07.856442",
    "2021-10-05 10:51:07.856444",
    "2021-10-05 10:51:07.856447",
    "2021-10-05 10:51:07.856450",
    "2021-10-05 10:51:07.856452",
    "2021-10-05 10:51:07.856455",
    "2021-10-05 10:51:07.856457",
    "2021-10-05 10:51:07.856460",
    "2021-10-05 10:51:07.856462",
    "2021-10-05 10:51:07.856465",
    "2021-10-05 10:51:07.856468",
    "2021-10-05 10:51:07.856471",
    "2021-10-05 10:51:07.856473",
    "2021-10-05 10:51:07.856476",
    "2021-10-05 10:51:07.856478",
    "2021-10-05 10:51:07.856481",
    "2021-10-05 10:51:07.856483",
    "2021-10-05 10:51:07.856486",
    "2021-10-05 10:51:07.856488",
    "2021-10-05 10:51:07.856491",
    "2021-10-05 10:51:07.856494",
    "2021-10-05 10:51:07.856497",
    "2021-10-05 10:51:07.856499",
    "2021-10-05 10:51:07.856502",
    "2021-10-05 10:51:07.856504",
    "2021-10-05 10:51:07.856507",
    "2021-10-05 10:51:07.856509",
    "2021-10-05 10:51:07.856512",
    "2021-10-05 10:51:07.856514",
    "2021-10-05 10:51:07.856517",
    "2021-10-05 10:51:07.856519",
    "2021-10-05 10:51:07.856522",
    "2021-10-05 10:51:07.856525",
    "2021-10-05 10:51:07.856528",
    "2021-10-05 10:51:07.856530",
    "2021-10-05 10:51:07.856533",
    "2021-10-05 10:51:07.856535",
    "2021-10-05 10:51:07.856538",
    "2021-10-05 10:51:07.856540",
    "2021-10-05 10:51:07.856543",
    "2021-10-05 10:51:07.856545",
    "2021-10-05 10:51:07.856548",
    "2021-10-05 10:51:07.856550",
    "2021-10-05 10:51:07.856554",
    "2021-10-05 10:51:07.856557",
    "2021-10-05 10:51:07.856560",
    "2021-10-05 10:51:07.856562",
    "2021-10-05 10:51:07.856565",
    "2021-10-05 10:51:07.856567",
    "2021-10-05 10:51:07.856570",
    "2021-10-05 10:51:07.856572",
    "2021-10-05 10:51:07.856575",
    "2021-10-05 10:51:07.856578",
    "2021-10-05 10:51:07.856580",
    "2021-10-05 10:51:07.856583",
    "2021-10-05 10:51:07.856587",
    "2021-10-05 10:51:07.856591",
    "2021-10-05 10:51:07.856593",
    "2021-10-05 10:51:07.856596",
    "2021-10-05 10:51:07.856599",
    "2021-10-05 10:51:07.856601",
    "2021-10-05 10:51:07.856604",
    "2021-10-05 10:51:07.856606",
    "2021-10-05 10:51:07.856609",
    "2021-10-05 10:51:07.856611",
    "2021-10-05 10:51:07.856614",
    "2021-10-05 10:51:07.856617",
    "2021-10-05 10:51:07.856619",
    "2021-10-05 10:51:07.856622",
    "2021-10-05 10:51:07.856624",
    "2021-10-05 10:51:07.856627",
    "2021-10-05 10:51:07.856631",
    "2021-10-05 10:51:07.856634",
    "2021-10-05 10:51:07.856636",
    "2021-10-05 10:51:07.856639",
    "2021-10-05 10:51:07.856641",
    "2021-10-05 10:51:07.856644",
    "2021-10-05 10:51:07.856646",
    "2021-10-05 10:51:07.856649",
    "2021-10-05 10:51:07.856652",
    "2021-10-05 10:51:07.856654",
    "2021-10-05 10:51:07.856657",
    "2021-10-05 10:51:07.856659",
    "2021-10-05 10:51:07.856662",
    "2021-10-05 10:51:07.856664",
    "2021-10-05 10:51:07.856667",
    "2021-10-05 10:51:07.856670",
    "2021-10-05 10:51:07.856674",
    "2021-10-05 10:51:07.856677",
    "2021-10-05 10:51:07.856680",
    "2021-10-05 10:51:07.856683",
    "2021-10-05 10:51:07.856685",
    "2021-10-05 10:51:07.856688",
    "2021-10-05 10:51:07.856691",
    "2021-10-05 10:51:07.856693",
    "2021-10-05 10:51:07.856696",
    "2021-10-05 10:51:07.856699",
    "2021-10-05 10:51:07.856701",
    "2021-10-05 10:51:07.856704",
    "2021-10-05 10:51:07.856706",
    "2021-10-05 10:51:07.856709",
    "2021-10-05 10:51:07.856711",
    "2021-10-05 10:51:07.856714",
    "2021-10-05 10:51:07.856717",
    "2021-10-05 10:51:07.856719",
    "2021-10-05 10:51:07.856722",
    "2021-10-05 10:51:07.856724",
    "2021-10-05 10:51:07.856729",
    "2021-10-05 10:51:07.856732",
    "2021-10-05 10:51:07.856735",
    "2021-10-05 10:51:07.856737",
    "2021-10-05 10:51:07.856740",
    "2021-10-05 10:51:07.856743",
    "2021-10-05 10:51:07.856745",
    "2021-10-05 10:51:07.856748",
    "2021-10-05 10:51:07.856751",
    "2021-10-05 10:51:07.856753",
    "2021-10-05 10:51:07.856756",
    "2021-10-05 10:51:07.856759",
    "2021-10-05 10:51:07.856761",
    "2021-10-05 10:51:07.856764",
    "2021-10-05 10:51:07.856766",
    "2021-10-05 10:51:07.856769",
    "2021-10-05 10:51:07.856772",
    "2021-10-05 10:51:07.856774",
    "2021-10-05 10:51:07.856777",
    "2021-10-05 10:51:07.856779",
    "2021-10-05 10:51:07.856782",
    "2021-10-05 10:51:07.856785",
    "2021-10-05 10:51:07.856788"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $timeValues[$i]
}
